# LOQ4254.xlsx syllabus update: fix the "Objetivos" body text, insert a new
# row for the professor name under "Docentes responsaveis:", and re-populate
# every content cell from "Programa resumido" through "Bibliografia" so each
# label lines back up with its real text (the source file had the values
# shifted by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ("Objetivos:") keeps its place but gets the real course objective text
$ws.Range("B10").Value = 'Conferir aos alunos uma visão geral da Indústria Siderúrgica, Metalúrgica, Mecânica e correlatas, bem como das principais características dos processos e arranjos produtivos destas indústrias.'
$ws.Range("C10").Value = 'Conferir aos alunos uma visão geral da Indústria Siderúrgica, Metalúrgica, Mecânica e correlatas, bem como das principais características dos processos e arranjos produtivos destas indústrias.'

# --- Insert a new row above the old row 13 ("Programa resumido:").
# Rows 13-23 shift down to 14-24; row heights/styles move with their content.
$ws.Rows.Item(13).Insert()

# Excel copies the formatting of the row above into the new row; row 13 should
# have no A cell at all (the "Docentes responsaveis:" label stays on row 12), so
# drop the stray cell the insert created.
$ws.Range("A13").Clear()

# B13/C13 are brand-new cells - copy the B/C cell formatting (styles 2/3) from
# another data row instead of re-describing fonts/alignment by hand, so we reuse
# the existing style entries rather than minting new ones.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 13 holds the professor name (label "Docentes responsaveis:" is row 12)
$ws.Range("B13").Value = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Range("C13").Value = '5840560 - Marco Antonio Carvalho Pereira'

# Rows 14, 16, 19, 20, 21, 22 already inherited the correct B/C styles from the
# shift caused by the insert above - only their text content needs fixing.

# Row 14: "Programa resumido:"
$ws.Range("B14").Value = '1) Processos da Indústria Siderúrgica. 2) Processos da Indústria Metalúrgica. 3) Processos da Indústria Mecânica.  4) Processos Industriais em Geral, exceto da Indústria Química'
$ws.Range("C14").Value = '1) Processos da Indústria Siderúrgica. 2) Processos da Indústria Metalúrgica. 3) Processos da Indústria Mecânica.  4) Processos Industriais em Geral, exceto da Indústria Química'

# Row 16: "Programa:"
$ws.Range("B16").Value = '1) Processos da Indústria Siderúrgica: Obtenção de Aços. Demais Processos.2) Processos da Indústria Metalúrgica: Processos de Fundição. Demais Processos.3) Processos da Indústria Mecânica: Processos de Conformação Plástica. Conformação por Corte de Usinagem. Demais Processos4) Processos Industriais em geral: Indústria da Construção Civil Indústria farmacêutica, Indústria Automobilística, dentre outras.'
$ws.Range("C16").Value = '1) Processos da Indústria Siderúrgica: Obtenção de Aços. Demais Processos.2) Processos da Indústria Metalúrgica: Processos de Fundição. Demais Processos.3) Processos da Indústria Mecânica: Processos de Conformação Plástica. Conformação por Corte de Usinagem. Demais Processos4) Processos Industriais em geral: Indústria da Construção Civil Indústria farmacêutica, Indústria Automobilística, dentre outras.'

# Row 19: "Metodo:"
$ws.Range("B19").Value = 'Aulas expositivas. Seminários. Palestras feiras por profissionais de indústrias. Trabalhos em grupo. Debates e palestras.'
$ws.Range("C19").Value = 'Aulas expositivas. Seminários. Palestras feiras por profissionais de indústrias. Trabalhos em grupo. Debates e palestras.'

# Row 20: "Criterio:"
$ws.Range("B20").Value = 'Média aritmética das atividades avaliativas realizadas.'
$ws.Range("C20").Value = 'Média aritmética das atividades avaliativas realizadas.'

# Row 21: "Norma de recuperacao:"
$ws.Range("B21").Value = 'Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.'
$ws.Range("C21").Value = 'Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.'

# Row 22: "Bibliografia:" (previously blank of content - now the reading list)
$ws.Range("B22").Value = '1. Marcelo Breda Mourão et al. Introdução à Siderurgia, ABM, São Paulo, 20072. Fathi Habashi. Extractive Metallurgy, Gordon and Breach Science Publishers, 1986. 3. Luiz Antônio de Araújo. Manual de siderurgia - produção, Editora Arte & Ciência, São Paulo, 1997. 4. Alan H. Cottrell. Introdução à metalurgia, 2a edição, Fundação Calouste Gulbenkian, Lisboa, 1975.5. ASM Handbook Vol. 15 Casting - 1988 , Foundry Technology P.R. Beeley, 19726. John Campbell. Casting Butterworth-Heinemann, 19917. M. Siegel, Fundição. ABM, S.Paulo, 1979. 8. Amauri Garcia. Solidificação: Fundamentos e Aplicações, Editora da UNICAMP, Campinas, SP, 20089. Mauricio Prates de Campos Filho e Graeme John Davies Solidificação e Fundição de Metais e suas Ligas, Livros Técnicos e Científicos, Rio de Janeiro.10. AVITZUR, B. Metal Forming: processes and analysis – TATA Mc Graw-Hill Publishing Company Limited; New Delhi, 1977.11. RODRIGUES, J.; MARTINS, P. Tecnologia Mecânica: Tecnologia da deformação plástica. Aplicações Industriais. Escolar Editora, v.1 e v.2, 2010.12. CETLIN, P.R.; HELMAN, H. Fundamentos da conformação mecânicas dos metais. Ed. Artliber Ltda, 260p., 2005.13. BRESCIANI FILHO, E.; ZAVAGLIA, C.A.C.; NERY, F.A.C.; BOTTON, S.T. Conformação plástica dos metais. Ed. Unicamp, v.1 e v.2, 1986.14. DINIZ, A.E.; MARCONDES, F.C.; COPPINI, N.L. Tecnologia da usinagem dos materiais. Ed. Artlebet Ltda., 244p., 2000.'
$ws.Range("C22").Value = '1. Marcelo Breda Mourão et al. Introdução à Siderurgia, ABM, São Paulo, 20072. Fathi Habashi. Extractive Metallurgy, Gordon and Breach Science Publishers, 1986. 3. Luiz Antônio de Araújo. Manual de siderurgia - produção, Editora Arte & Ciência, São Paulo, 1997. 4. Alan H. Cottrell. Introdução à metalurgia, 2a edição, Fundação Calouste Gulbenkian, Lisboa, 1975.5. ASM Handbook Vol. 15 Casting - 1988 , Foundry Technology P.R. Beeley, 19726. John Campbell. Casting Butterworth-Heinemann, 19917. M. Siegel, Fundição. ABM, S.Paulo, 1979. 8. Amauri Garcia. Solidificação: Fundamentos e Aplicações, Editora da UNICAMP, Campinas, SP, 20089. Mauricio Prates de Campos Filho e Graeme John Davies Solidificação e Fundição de Metais e suas Ligas, Livros Técnicos e Científicos, Rio de Janeiro.10. AVITZUR, B. Metal Forming: processes and analysis – TATA Mc Graw-Hill Publishing Company Limited; New Delhi, 1977.11. RODRIGUES, J.; MARTINS, P. Tecnologia Mecânica: Tecnologia da deformação plástica. Aplicações Industriais. Escolar Editora, v.1 e v.2, 2010.12. CETLIN, P.R.; HELMAN, H. Fundamentos da conformação mecânicas dos metais. Ed. Artliber Ltda, 260p., 2005.13. BRESCIANI FILHO, E.; ZAVAGLIA, C.A.C.; NERY, F.A.C.; BOTTON, S.T. Conformação plástica dos metais. Ed. Unicamp, v.1 e v.2, 1986.14. DINIZ, A.E.; MARCONDES, F.C.; COPPINI, N.L. Tecnologia da usinagem dos materiais. Ed. Artlebet Ltda., 244p., 2000.'

